# edit.ps1
# Refresh the "想去人数" (interested-count) figures in column F of the
# "展览" (Exhibitions) and "全部类型" (All types) worksheets, matching the
# data snapshot committed as "Update gh-pages to output generated at 456a3b4".
# Only the numeric values in column F change; nothing else in the workbook
# is touched.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Update sheet '展览' (Exhibitions) - column F (想去人数 / interested-count)
$ws1.Cells.Item(2, 6).Value = 830
$ws1.Cells.Item(3, 6).Value = 565
$ws1.Cells.Item(4, 6).Value = 313
$ws1.Cells.Item(11, 6).Value = 1214
$ws1.Cells.Item(14, 6).Value = 904
$ws1.Cells.Item(15, 6).Value = 888
$ws1.Cells.Item(20, 6).Value = 797
$ws1.Cells.Item(21, 6).Value = 1748
$ws1.Cells.Item(22, 6).Value = 3050
$ws1.Cells.Item(23, 6).Value = 893
$ws1.Cells.Item(24, 6).Value = 88
$ws1.Cells.Item(25, 6).Value = 2267
$ws1.Cells.Item(27, 6).Value = 2
$ws1.Cells.Item(28, 6).Value = 3123
$ws1.Cells.Item(29, 6).Value = 634
$ws1.Cells.Item(30, 6).Value = 624
$ws1.Cells.Item(31, 6).Value = 18
$ws1.Cells.Item(32, 6).Value = 92
$ws1.Cells.Item(33, 6).Value = 741
$ws1.Cells.Item(35, 6).Value = 137
$ws1.Cells.Item(36, 6).Value = 13
$ws1.Cells.Item(37, 6).Value = 103
$ws1.Cells.Item(38, 6).Value = 1114
$ws1.Cells.Item(39, 6).Value = 1802
$ws1.Cells.Item(40, 6).Value = 405
$ws1.Cells.Item(43, 6).Value = 201
$ws1.Cells.Item(44, 6).Value = 136
$ws1.Cells.Item(45, 6).Value = 185

# Update sheet '全部类型' (All types) - column F (想去人数 / interested-count)
$ws4.Cells.Item(2, 6).Value = 830
$ws4.Cells.Item(3, 6).Value = 565
$ws4.Cells.Item(4, 6).Value = 313
$ws4.Cells.Item(10, 6).Value = 1214
$ws4.Cells.Item(12, 6).Value = 904
$ws4.Cells.Item(13, 6).Value = 888
$ws4.Cells.Item(19, 6).Value = 797
$ws4.Cells.Item(20, 6).Value = 1748
$ws4.Cells.Item(21, 6).Value = 3050
$ws4.Cells.Item(22, 6).Value = 893
$ws4.Cells.Item(23, 6).Value = 88
$ws4.Cells.Item(25, 6).Value = 2267
$ws4.Cells.Item(26, 6).Value = 3123
$ws4.Cells.Item(27, 6).Value = 634
$ws4.Cells.Item(28, 6).Value = 624
$ws4.Cells.Item(30, 6).Value = 18
$ws4.Cells.Item(34, 6).Value = 92
$ws4.Cells.Item(36, 6).Value = 741
$ws4.Cells.Item(38, 6).Value = 137
$ws4.Cells.Item(39, 6).Value = 103
$ws4.Cells.Item(41, 6).Value = 1114
$ws4.Cells.Item(42, 6).Value = 1802
$ws4.Cells.Item(44, 6).Value = 405
$ws4.Cells.Item(46, 6).Value = 201
$ws4.Cells.Item(47, 6).Value = 136
$ws4.Cells.Item(48, 6).Value = 185
